$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/1996 -> 01/01/2022
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"

# Docentes responsáveis (second docente) change
$ws.Range("B14").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C14").Value = "519033 - Carlos Yujiro Shigue"

# Método:
$ws.Range("B20").Value = "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes."
$ws.Range("C20").Value = "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes."

# Critério:
$ws.Range("B21").Value = "A média do semestre será computada com base na relação: M=(A1+A2)/2"
$ws.Range("C21").Value = "A média do semestre será computada com base na relação: M=(A1+A2)/2"

# Norma de recuperação:
$ws.Range("B22").Value = "Não cabe recuperação."
$ws.Range("C22").Value = "Não cabe recuperação."
